# Generate Report for Handback
#
# - Flips every "Ready for handoff" status cell to "Handed back: in sync with en-US"
# - Stamps the "Latest Handback DateTime" for the zh-cn and de-de rows
# - Adds the "Latest Target File" (F) and "Latest Handback File" (G) columns, with
#   hyperlinks mirroring the existing "Source File Name" (A) / "Latest Handoff File" (D)
#   links, for both data rows of the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

function Get-HyperlinkAddress($sheet, $cellAddr) {
    $target = '$' + $cellAddr.Substring(0,1) + '$' + $cellAddr.Substring(1)
    foreach ($hl in $sheet.Hyperlinks) {
        if ($hl.Range.Address() -eq $target) {
            return $hl.Address
        }
    }
    return $null
}

# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US" everywhere.
$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

foreach ($sheetName in @("Overview", "zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $used = $ws.UsedRange
    $rows = $used.Row + $used.Rows.Count - 1
    $cols = $used.Column + $used.Columns.Count - 1
    for ($r = $used.Row; $r -le $rows; $r++) {
        for ($c = $used.Column; $c -le $cols; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            if ($cell.Value2 -eq $oldStatus) {
                $cell.Value = $newStatus
            }
        }
    }
}

# 2) Latest Handback DateTime (column H) for the data rows.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsZh.Range("H2").Value = "2016-03-21 00:54:12"
$wsZh.Range("H3").Value = "2016-03-21 00:54:12"

$wsDe.Range("H2").Value = "2016-03-21 00:54:18"
$wsDe.Range("H3").Value = "2016-03-21 00:54:18"

# 3) New "Latest Target File" (F) / "Latest Handback File" (G) columns.
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $urlA2 = Get-HyperlinkAddress $ws "A2"
    $urlD2 = Get-HyperlinkAddress $ws "D2"

    $textF = $ws.Range("A2").Value2
    $textG = $ws.Range("D2").Value2

    foreach ($r in @(2, 3)) {
        $fCell = $ws.Cells.Item($r, 6)
        $gCell = $ws.Cells.Item($r, 7)

        $ws.Hyperlinks.Add($fCell, $urlA2, "", "", $textF) | Out-Null
        $ws.Hyperlinks.Add($gCell, $urlD2, "", "", $textG) | Out-Null
    }
}
